# Atualiza dados da BIBI - metricas_retencao_anual
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Linha 22: cohort 2020, period_index 5 -> num_customers 30 -> 31
$ws.Range("C22").Value = 31
$ws.Range("E22").Value = 31 / 2654

# Linha 36: cohort 2024, period_index 1 -> num_customers 136 -> 138
$ws.Range("C36").Value = 138
$ws.Range("E36").Value = 138 / 1930

# Linha 37: cohort 2024, period_index 0 -> num_customers 870 -> 877, cohort_size 870 -> 877
$ws.Range("C37").Value = 877
$ws.Range("D37").Value = 877
$ws.Range("E37").Value = 877 / 877
